$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Sheet2: fill in the bill data ---
$ws2.Range("A1").Value = 18
$ws2.Range("D1").Value = 279

$ws2.Range("A2").Value = 44
$ws2.Range("D2").Value = 69

$ws2.Range("A3").Value = 46
$ws2.Range("D3").Value = 208

$ws2.Range("A4").Value = 18
$ws2.Range("D4").Value = 110

$ws2.Range("A5").Value = 22
$ws2.Range("D5").Value = 598

$ws2.Range("A6").Formula = "=SUM(A1:A5)"
$ws2.Range("B6").Value = "RMB"
$ws2.Range("D6").Value = 159

$ws2.Range("D7").Value = 159

$ws2.Range("D8").Value = 319
$ws2.Range("G8").Value = 42
$ws2.Range("H8").Value = "hkd"

$ws2.Range("D9").Value = 3888

$ws2.Range("D10").Formula = "=SUM(D1:D9)"
$ws2.Range("E10").Value = "hkd"
$ws2.Range("M10").Value = 5789

$ws2.Range("M11").Value = 668

$ws2.Range("M12").Value = 42

$ws2.Range("E13").Value = 168
$ws2.Range("M13").Formula = "=SUM(M10:M12)"

$ws2.Range("E14").Value = 500

# --- Sheet1: move selection off the old cell ---
$ws1.Range("C1:C11").Select() | Out-Null

# --- Sheet2 becomes the active/visible sheet with its own selection ---
$ws2.Activate() | Out-Null
$ws2.Range("K27").Select() | Out-Null
